# Update gh-pages to output generated at 456a3b4
# Updates "想去人数" (F column) figures on the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 11779
$ws1.Range("F5").Value = 345
$ws1.Range("F7").Value = 11716
$ws1.Range("F9").Value = 1168
$ws1.Range("F11").Value = 32
$ws1.Range("F12").Value = 1770
$ws1.Range("F13").Value = 5803

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 11779
$ws4.Range("F7").Value = 345
$ws4.Range("F9").Value = 11716
$ws4.Range("F11").Value = 1168
$ws4.Range("F13").Value = 32
$ws4.Range("F14").Value = 1770
$ws4.Range("F16").Value = 5803
